# Apply "Natmi following Dr Hou advice" data update to sheet1.
# Target shape: rows 2-13 (12 data rows), columns A-T (20 columns).
# Column layout: A=Sending cluster, B=Ligand symbol, C=Receptor symbol,
# D=Target cluster, E..T = numeric metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "F2", "F2r", "ECs", 3, 1, 0.6740386666666667, 2.022116, 0.3058319194585966, 0.3058319194585966, 3, 1, 8.607073666666667, 25.821221, 0.09431423806264486, 0.09431423806264484, 5.801500458181778, 52.213504123636, 0.02884430445897371, 0.0288443044589737),
  @("ECs", "F2", "F2r", "FAPs", 3, 1, 0.6740386666666667, 2.022116, 0.3058319194585966, 0.3058319194585966, 3, 1, 71.03134166666666, 213.094025, 0.7783443162342011, 0.778344316234201, 47.87787082854445, 430.9008374569, 0.2380425362335946, 0.2380425362335946),
  @("ECs", "F2", "F2r", "sCs", 3, 1, 0.6740386666666667, 2.022116, 0.3058319194585966, 0.3058319194585966, 3, 1, 11.62112133333333, 34.863364, 0.1273414457031541, 0.1273414457031541, 7.833085128691556, 70.497766158224, 0.03894507876602828, 0.03894507876602828),
  @("FAPs", "F2", "F2r", "ECs", 3, 1, 0.9731926666666667, 2.919578, 0.4415672215387696, 0.4415672215387696, 3, 1, 8.607073666666667, 25.821221, 0.09431423806264486, 0.09431423806264484, 8.376340973859778, 75.38706876473799, 0.04164607605286816, 0.04164607605286815),
  @("FAPs", "F2", "F2r", "FAPs", 3, 1, 0.9731926666666667, 2.919578, 0.4415672215387696, 0.4415672215387696, 3, 1, 71.03134166666666, 213.094025, 0.7783443162342011, 0.778344316234201, 69.12718081349445, 622.14462732145, 0.3436913371200296, 0.3436913371200296),
  @("FAPs", "F2", "F2r", "sCs", 3, 1, 0.9731926666666667, 2.919578, 0.4415672215387696, 0.4415672215387696, 3, 1, 11.62112133333333, 34.863364, 0.1273414457031541, 0.1273414457031541, 11.30959006004356, 101.786310540392, 0.05622980836587184, 0.05622980836587184),
  @("M2", "F2", "F2r", "ECs", 1, 0.3333333333333333, 0.048595, 0.145785, 0.02204903496054208, 0.02204903496054208, 3, 1, 8.607073666666667, 25.821221, 0.09431423806264486, 0.09431423806264484, 0.4182607448316666, 3.764346703485, 0.002079537932320145, 0.002079537932320144),
  @("M2", "F2", "F2r", "FAPs", 1, 0.3333333333333333, 0.048595, 0.145785, 0.02204903496054208, 0.02204903496054208, 3, 1, 71.03134166666666, 213.094025, 0.7783443162342011, 0.778344316234201, 3.451768048291667, 31.065912434625, 0.01716174103998712, 0.01716174103998712),
  @("M2", "F2", "F2r", "sCs", 1, 0.3333333333333333, 0.048595, 0.145785, 0.02204903496054208, 0.02204903496054208, 3, 1, 11.62112133333333, 34.863364, 0.1273414457031541, 0.1273414457031541, 0.5647283911933333, 5.08255552074, 0.002807755988234815, 0.002807755988234815),
  @("sCs", "F2", "F2r", "ECs", 3, 1, 0.508125, 1.524375, 0.2305518240420917, 0.2305518240420917, 3, 1, 8.607073666666667, 25.821221, 0.09431423806264486, 0.09431423806264484, 4.373469306875, 39.361223761875, 0.02174431961848284, 0.02174431961848284),
  @("sCs", "F2", "F2r", "FAPs", 3, 1, 0.508125, 1.524375, 0.2305518240420917, 0.2305518240420917, 3, 1, 71.03134166666666, 213.094025, 0.7783443162342011, 0.778344316234201, 36.092800484375, 324.835204359375, 0.1794487018405897, 0.1794487018405897),
  @("sCs", "F2", "F2r", "sCs", 3, 1, 0.508125, 1.524375, 0.2305518240420917, 0.2305518240420917, 3, 1, 11.62112133333333, 34.863364, 0.1273414457031541, 0.1273414457031541, 5.9049822775, 53.1448404975, 0.02935880258301915, 0.02935880258301915)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = $startRow + $i
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $row[$j]
  }
}
